$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Get-MCASActivity" (row 24, General category) is now functional: give it the
# same yellow "done" row formatting used by other completed rows (e.g. row 26)
# and mark it with an "x" in the same columns that row uses.
$ws.Range("A26:Q26").Copy()
$ws.Range("A24:Q24").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C24:J24").Value = "x"
$ws.Range("N24").Value = "x"
$ws.Range("Q24").Value = "x"

# "ConvertFrom-MCASTimestamp" (row 26) picks up a few more functional markers.
$ws.Range("N26").Value = "x"
$ws.Range("O26").Value = "x"
$ws.Range("Q26").Value = "x"

# The author's selection ended up on Q26 after this edit.
$ws.Range("Q26").Select()
